$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "53.337.58"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -11.95%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.319.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -20.04%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "438.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -17.08%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "121.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -15.29%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.995"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.44%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.473"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -14.65%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.312.38"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -20.48%  "

$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0915"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -15.15%  "

$ws.Range("B11").Value = "Toncoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.27"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -12.83%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.308"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -14.76%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.122"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.55%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.709.50"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -20.50%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "53.357.09"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -11.91%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.08"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -16.30%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000120"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -14.94%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.327.84"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -19.96%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.97"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -21.19%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "301.60"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -16.60%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.12"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -22.14%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.17%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.64"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.84%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -19.02%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "55.38"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -14.28%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.995"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.30%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.153"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -14.30%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.368"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -19.17%  "

$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.97"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -11.26%  "

$ws.Range("B30").Value = "USDe"
$ws.Range("C30").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.995"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.46%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0₃0697"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -18.22%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "144.02"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.76%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "17.16"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -13.09%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.35"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -19.71%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.67"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -16.10%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.53"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -19.14%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.824"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -17.70%  "

$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -16.71%  "

$ws.Range("B39").Value = "FirstDigitalUSD"
$ws.Range("C39").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.994"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.31%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "32.89"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -12.96%  "

$ws.Range("E41").Value = "  -0.59%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.15"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -15.36%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.924.48"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -16.17%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.21"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -18.29%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0497"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -14.61%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.520"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -19.79%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0209"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -11.99%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0829"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -10.09%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "15.67"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -23.22%  "

$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.95"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -20.36%  "

$ws.Range("B51").Value = "ZEEBU"
$ws.Range("C51").Value = "https://coinranking.com/coin/B5-YKN_zB+zeebu-zbu"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.63"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.19%  "
